# Apply cryptos list update (values and two row swaps) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.710.98"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").Value = "2.701.90"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").Value = "2.726.00"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.130"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").Value = "3.175.16"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "60.672.02"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "2.719.43"
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "348.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  +5.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.942"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.89%  "
$ws.Range("E39").Value = "  +7.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.613"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.148.61"
$ws.Range("E46").Value = "  +7.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.59%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.69%  "
